$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Broken Links" block (rows 58-61), mirroring the existing Links section rows.
$ws.Range("A58").Value = "brokenLinks"

$ws.Range("A59").Value = "brokenLinksNav"
$ws.Range("B59").Value = '//*[@id="item-6"]'
$ws.Range("C59").Value = "By.xpath"

$ws.Range("A60").Value = "brokenLink"
$ws.Range("B60").Value = '//*[@id="app"]/div/div/div/div[2]/div[2]/a[2]'
$ws.Range("C60").Value = "By.xpath"

$ws.Range("A61").Value = "brokenLinkScroll"
$ws.Range("B61").Value = '//*[@id="app"]/div/div/div/div[2]/div[2]/h1'
$ws.Range("C61").Value = "By.xpath"

# Selection state as left by the author after making the edit.
$ws.Range("A55").Select() | Out-Null
